$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from H1 onto the two new header cells, then set values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-24: I = 1, J = same value as H
$hValues = @{
    2 = 5
    3 = 5
    4 = 5
    5 = 6
    6 = 5
    7 = 7
    8 = 5
    9 = 5
    10 = 6
    11 = 8
    12 = 5
    13 = 7
    14 = 6
    15 = 5
    16 = 7
    17 = 7
    18 = 5
    19 = 5
    20 = 5
    21 = 5
    22 = 5
    23 = 4
    24 = 4
}

foreach ($r in $hValues.Keys) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hValues[$r]
}

# Row 25 is a special case (does not follow the I=1,J=H pattern)
$ws.Cells.Item(25, 9).Value = 3
$ws.Cells.Item(25, 10).Value = 4
